$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2679.25
$ws.Range("J32").Value = 2389.1667
$ws.Range("L32").Value = 2389.1667
$ws.Range("N32").Value = -3041.1667
$ws.Range("H40").Value = 3605.12
$ws.Range("I40").Value = 4980.154
$ws.Range("J40").Value = 2115.5
$ws.Range("K40").Value = 4980.154
$ws.Range("L40").Value = 2115.5
$ws.Range("M40").Value = -4805.154
$ws.Range("N40").Value = -2465.5
$ws.Range("H57").Value = 47000
$ws.Range("J57").Value = 47000
$ws.Range("L57").Value = 141000
$ws.Range("N57").Value = -141998
$ws.Range("H76").Value = 6504.5454
$ws.Range("J76").Value = 7499.7144
$ws.Range("L76").Value = 7499.7144
$ws.Range("N76").Value = -8129.7144
$ws.Range("H79").Value = 6504.5454
$ws.Range("J79").Value = 7499.7144
$ws.Range("L79").Value = 7499.7144
$ws.Range("N79").Value = -9683.714400000001
$ws.Range("H94").Value = 6864.6665
$ws.Range("I94").Value = 2943.2727
$ws.Range("K94").Value = 2943.2727
$ws.Range("M94").Value = -2492.2727
$ws.Range("H100").Value = 3814.577
$ws.Range("I100").Value = 1590.1875
$ws.Range("J100").Value = 7373.6
$ws.Range("K100").Value = 1590.1875
$ws.Range("L100").Value = 7373.6
$ws.Range("M100").Value = -1049.1875
$ws.Range("N100").Value = -8455.6
$ws.Range("H132").Value = 66608.78
$ws.Range("I132").Value = 74249.28
$ws.Range("K132").Value = 222747.84
$ws.Range("M132").Value = -220217.84

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4001.6667
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3788
$ws.Range("H132").Value = 4039.1428
$ws.Range("I132").Value = 3057.625
$ws.Range("K132").Value = 9172.875
$ws.Range("M132").Value = -6642.875
$ws.Range("H136").Value = 4001.6667
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 1024.4117
$ws.Range("I86").Value = 887.4
$ws.Range("K86").Value = 887.4
$ws.Range("M86").Value = 235.6
$ws.Range("H89").Value = 1024.4117
$ws.Range("I89").Value = 887.4
$ws.Range("K89").Value = 4437
$ws.Range("M89").Value = 1179
$ws.Range("H107").Value = 2139.1843
$ws.Range("J107").Value = 3375.3333
$ws.Range("L107").Value = 3375.3333
$ws.Range("N107").Value = -7215.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 69081
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 69081
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 3071.6667
$ws.Range("I132").Value = 2355.1428
$ws.Range("J132").Value = 4074.8
$ws.Range("K132").Value = 7065.428400000001
$ws.Range("L132").Value = 12224.4
$ws.Range("M132").Value = -4535.428400000001
$ws.Range("N132").Value = -17284.4
$ws.Range("H139").Value = 82677.25
$ws.Range("I139").Value = 82677.25
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 82677.25
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -77537.25
$ws.Range("N139").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2882.75
$ws.Range("J129").Value = 3396.1538
$ws.Range("L129").Value = 10188.4614
$ws.Range("N129").Value = -20188.4614
$ws.Range("H141").Value = 2641.3333
$ws.Range("I141").Value = 2641.3333
$ws.Range("K141").Value = 7923.999899999999
$ws.Range("M141").Value = -2743.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3577.64
$ws.Range("I70").Value = 3459.4375
$ws.Range("K70").Value = 3459.4375
$ws.Range("M70").Value = -3189.4375
$ws.Range("H73").Value = 3577.64
$ws.Range("I73").Value = 3459.4375
$ws.Range("K73").Value = 3459.4375
$ws.Range("M73").Value = -2523.4375
$ws.Range("H80").Value = 192463.11
$ws.Range("I80").Value = 304158.53
$ws.Range("J80").Value = 2580.9
$ws.Range("K80").Value = 304158.53
$ws.Range("L80").Value = 2580.9
$ws.Range("M80").Value = -303160.53
$ws.Range("N80").Value = -4576.9
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 192463.11
$ws.Range("I83").Value = 304158.53
$ws.Range("J83").Value = 2580.9
$ws.Range("K83").Value = 1520792.65
$ws.Range("L83").Value = 12904.5
$ws.Range("M83").Value = -1515800.65
$ws.Range("N83").Value = -22888.5
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1700.359
$ws.Range("I40").Value = 1800.9117
$ws.Range("J40").Value = 1016.6
$ws.Range("K40").Value = 1800.9117
$ws.Range("L40").Value = 1016.6
$ws.Range("M40").Value = -1664.9117
$ws.Range("N40").Value = -1288.6
$ws.Range("H76").Value = 10288
$ws.Range("J76").Value = 10288
$ws.Range("L76").Value = 10288
$ws.Range("N76").Value = -10964
$ws.Range("H79").Value = 10288
$ws.Range("J79").Value = 10288
$ws.Range("L79").Value = 10288
$ws.Range("N79").Value = -12628
$ws.Range("H93").Value = 4666.3335
$ws.Range("I93").Value = 4499.5
$ws.Range("K93").Value = 4499.5
$ws.Range("M93").Value = -3251.5
$ws.Range("H100").Value = 22400.8
$ws.Range("H122").Value = 3998.756
$ws.Range("I122").Value = 3774.5405
$ws.Range("K122").Value = 11323.6215
$ws.Range("M122").Value = -8873.621500000001
$ws.Range("H136").Value = 3019.6843
$ws.Range("I136").Value = 2391.6667
$ws.Range("J136").Value = 5374.75
$ws.Range("K136").Value = 7175.000100000001
$ws.Range("L136").Value = 16124.25
$ws.Range("M136").Value = -4625.000100000001
$ws.Range("N136").Value = -21224.25
$ws.Range("H137").Value = 89999.5
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4937.778
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H136").Value = 7575
$ws.Range("I136").Value = 7206.25
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 21618.75
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -19068.75
$ws.Range("N136").Value = -41100
